$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '66.351.32'
$c.Style = "Normal"
$ws.Range("E2").Value = '  -0.76%  '
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '2.572.88'
$c.Style = "Normal"
$ws.Range("E3").Value = '  -1.81%  '
$ws.Range("E4").Value = '  -0.03%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '580.84'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -1.37%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '164.97'
$c.Style = "Normal"
$ws.Range("E6").Value = '  +0.14%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("E8").Value = '  -1.62%  '
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '2.572.07'
$c.Style = "Normal"
$ws.Range("E9").Value = '  -1.81%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '0.137'
$c.Style = "Normal"
$ws.Range("E10").Value = '  -3.74%  '
$ws.Range("E11").Value = '  +0.29%  '
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '0.358'
$c.Style = "Normal"
$ws.Range("E12").Value = '  -0.65%  '
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '5.15'
$c.Style = "Normal"
$ws.Range("E13").Value = '  -1.19%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '26.75'
$c.Style = "Normal"
$ws.Range("E14").Value = '  -2.50%  '
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '3.040.98'
$c.Style = "Normal"
$ws.Range("E15").Value = '  -2.45%  '
$ws.Range("E16").Value = '  -1.99%  '
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '66.266.83'
$c.Style = "Normal"
$ws.Range("E17").Value = '  -0.87%  '
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '2.580.91'
$c.Style = "Normal"
$ws.Range("E18").Value = '  -2.04%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '11.43'
$c.Style = "Normal"
$ws.Range("E19").Value = '  -3.90%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '7.72'
$c.Style = "Normal"
$ws.Range("E20").Value = '  -3.94%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '350.22'
$c.Style = "Normal"
$ws.Range("E21").Value = '  -1.96%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '4.22'
$c.Style = "Normal"
$ws.Range("E22").Value = '  -2.37%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '4.58'
$c.Style = "Normal"
$ws.Range("E23").Value = '  -2.34%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '10.09'
$c.Style = "Normal"
$ws.Range("E25").Value = '  -7.63%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '1.88'
$c.Style = "Normal"
$ws.Range("E26").Value = '  -3.89%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '68.87'
$c.Style = "Normal"
$ws.Range("E27").Value = '  -2.43%  '
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '2.706.74'
$c.Style = "Normal"
$ws.Range("E28").Value = '  -1.77%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Range("E29").Value = '  +0.09%  '
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '0.0₃0982'
$c.Style = "Normal"
$ws.Range("E30").Value = '  -2.39%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '532.30'
$c.Style = "Normal"
$ws.Range("E31").Value = '  -3.43%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '7.96'
$c.Style = "Normal"
$ws.Range("E32").Value = '  +0.66%  '
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '1.32'
$c.Style = "Normal"
$ws.Range("E33").Value = '  -2.73%  '
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '1.84'
$c.Style = "Normal"
$ws.Range("E34").Value = '  -2.82%  '
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '0.131'
$c.Style = "Normal"
$ws.Range("E35").Value = '  -0.71%  '
$ws.Range("E36").Value = '  +0.08%  '
$ws.Range("E37").Value = '  -3.80%  '
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '156.82'
$c.Style = "Normal"
$ws.Range("E38").Value = '  -0.23%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '18.69'
$c.Style = "Normal"
$ws.Range("E39").Value = '  -2.13%  '
$ws.Range("E40").Value = '  -1.61%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '18.22'
$c.Style = "Normal"
$ws.Range("E41").Value = '  +1.95%  '
$ws.Range("E42").Value = '  -0.86%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '5.08'
$c.Style = "Normal"
$ws.Range("E43").Value = '  -2.16%  '
$ws.Range("E44").Value = '  +0.09%  '
$ws.Range("E45").Value = '  -3.77%  '
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '0.0₆0286'
$c.Style = "Normal"
$ws.Range("E46").Value = '  -2.79%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '148.46'
$c.Style = "Normal"
$ws.Range("E47").Value = '  -2.10%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '0.564'
$c.Style = "Normal"
$ws.Range("E48").Value = '  -3.47%  '
$ws.Range("E49").Value = '  -2.04%  '
$ws.Range("E50").Value = '  -1.73%  '
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '0.0759'
$c.Style = "Normal"
$ws.Range("E51").Value = '  -1.85%  '
